# Auto-generated Excel COM-interop script to apply cryptos.xlsx diff
# Commit: Updated symbol list on Mon Feb  6 00:43:10 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'329.89"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.37%"
$ws.Range("E2").Style = "Normal"
$ws.Range("F2").Value = "'6-2-2023"
$ws.Range("F2").Style = "Normal"
$ws.Range("G2").Value = "'0"
$ws.Range("G2").Style = "Normal"
$ws.Range("D3").Value = "'45.32"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'10.33%"
$ws.Range("E3").Style = "Normal"
$ws.Range("F3").Value = "'6-2-2023"
$ws.Range("F3").Style = "Normal"
$ws.Range("G3").Value = "'0"
$ws.Range("G3").Style = "Normal"
$ws.Range("D4").Value = "'5.566"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-2.09%"
$ws.Range("E4").Style = "Normal"
$ws.Range("F4").Value = "'6-2-2023"
$ws.Range("F4").Style = "Normal"
$ws.Range("G4").Value = "'0"
$ws.Range("G4").Style = "Normal"
$ws.Range("D5").Value = "'0.08113"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-2.94%"
$ws.Range("E5").Style = "Normal"
$ws.Range("F5").Value = "'6-2-2023"
$ws.Range("F5").Style = "Normal"
$ws.Range("G5").Value = "'0"
$ws.Range("G5").Style = "Normal"
$ws.Range("D6").Value = "'8.683"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-1.54%"
$ws.Range("E6").Style = "Normal"
$ws.Range("F6").Value = "'6-2-2023"
$ws.Range("F6").Style = "Normal"
$ws.Range("G6").Value = "'0"
$ws.Range("G6").Style = "Normal"
$ws.Range("D7").Value = "'1.917"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-3.80%"
$ws.Range("E7").Style = "Normal"
$ws.Range("F7").Value = "'6-2-2023"
$ws.Range("F7").Style = "Normal"
$ws.Range("G7").Value = "'0"
$ws.Range("G7").Style = "Normal"
$ws.Range("D8").Value = "'4.300"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-4.51%"
$ws.Range("E8").Style = "Normal"
$ws.Range("F8").Value = "'6-2-2023"
$ws.Range("F8").Style = "Normal"
$ws.Range("G8").Value = "'0"
$ws.Range("G8").Style = "Normal"
$ws.Range("D9").Value = "'2.746"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-5.65%"
$ws.Range("E9").Style = "Normal"
$ws.Range("F9").Value = "'6-2-2023"
$ws.Range("F9").Style = "Normal"
$ws.Range("G9").Value = "'0"
$ws.Range("G9").Style = "Normal"
$ws.Range("D10").Value = "'0.9477"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'2.69%"
$ws.Range("E10").Style = "Normal"
$ws.Range("F10").Value = "'6-2-2023"
$ws.Range("F10").Style = "Normal"
$ws.Range("G10").Value = "'0"
$ws.Range("G10").Style = "Normal"
$ws.Range("D11").Value = "'0.1196"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-4.06%"
$ws.Range("E11").Style = "Normal"
$ws.Range("F11").Value = "'6-2-2023"
$ws.Range("F11").Style = "Normal"
$ws.Range("G11").Value = "'0"
$ws.Range("G11").Style = "Normal"
$ws.Range("D12").Value = "'0.1904"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-2.65%"
$ws.Range("E12").Style = "Normal"
$ws.Range("F12").Value = "'6-2-2023"
$ws.Range("F12").Style = "Normal"
$ws.Range("G12").Value = "'0"
$ws.Range("G12").Style = "Normal"
$ws.Range("D13").Value = "'0.09821"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'5.19%"
$ws.Range("E13").Style = "Normal"
$ws.Range("F13").Value = "'6-2-2023"
$ws.Range("F13").Style = "Normal"
$ws.Range("G13").Value = "'0"
$ws.Range("G13").Style = "Normal"
$ws.Range("D14").Value = "'0.04119"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'4.22%"
$ws.Range("E14").Style = "Normal"
$ws.Range("F14").Value = "'6-2-2023"
$ws.Range("F14").Style = "Normal"
$ws.Range("G14").Value = "'0"
$ws.Range("G14").Style = "Normal"
$ws.Range("D15").Value = "'0.1067"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.30%"
$ws.Range("E15").Style = "Normal"
$ws.Range("F15").Value = "'6-2-2023"
$ws.Range("F15").Style = "Normal"
$ws.Range("G15").Value = "'0"
$ws.Range("G15").Style = "Normal"
$ws.Range("D16").Value = "'0.001280"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-0.92%"
$ws.Range("E16").Style = "Normal"
$ws.Range("F16").Value = "'6-2-2023"
$ws.Range("F16").Style = "Normal"
$ws.Range("G16").Value = "'0"
$ws.Range("G16").Style = "Normal"
$ws.Range("D17").Value = "'0.005930"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-2.81%"
$ws.Range("E17").Style = "Normal"
$ws.Range("F17").Value = "'6-2-2023"
$ws.Range("F17").Style = "Normal"
$ws.Range("G17").Value = "'0"
$ws.Range("G17").Style = "Normal"
$ws.Range("D18").Value = "'3.584"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'3.78%"
$ws.Range("E18").Style = "Normal"
$ws.Range("F18").Value = "'6-2-2023"
$ws.Range("F18").Style = "Normal"
$ws.Range("G18").Value = "'0"
$ws.Range("G18").Style = "Normal"
$ws.Range("D19").Value = "'0.3484"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-0.73%"
$ws.Range("E19").Style = "Normal"
$ws.Range("F19").Value = "'6-2-2023"
$ws.Range("F19").Style = "Normal"
$ws.Range("G19").Value = "'0"
$ws.Range("G19").Style = "Normal"
$ws.Range("D20").Value = "'8.621"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-4.53%"
$ws.Range("E20").Style = "Normal"
$ws.Range("F20").Value = "'6-2-2023"
$ws.Range("F20").Style = "Normal"
$ws.Range("G20").Value = "'0"
$ws.Range("G20").Style = "Normal"
$ws.Range("D21").Value = "'0.1364"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-0.93%"
$ws.Range("E21").Style = "Normal"
$ws.Range("F21").Value = "'6-2-2023"
$ws.Range("F21").Style = "Normal"
$ws.Range("G21").Value = "'0"
$ws.Range("G21").Style = "Normal"
$ws.Range("D22").Value = "'0.2588"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-1.91%"
$ws.Range("E22").Style = "Normal"
$ws.Range("F22").Value = "'6-2-2023"
$ws.Range("F22").Style = "Normal"
$ws.Range("G22").Value = "'0"
$ws.Range("G22").Style = "Normal"
$ws.Range("D23").Value = "'0.04358"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-1.52%"
$ws.Range("E23").Style = "Normal"
$ws.Range("F23").Value = "'6-2-2023"
$ws.Range("F23").Style = "Normal"
$ws.Range("G23").Value = "'0"
$ws.Range("G23").Style = "Normal"
$ws.Range("E24").Value = "'-0.69%"
$ws.Range("E24").Style = "Normal"
$ws.Range("F24").Value = "'6-2-2023"
$ws.Range("F24").Style = "Normal"
$ws.Range("G24").Value = "'0"
$ws.Range("G24").Style = "Normal"
$ws.Range("D25").Value = "'0.004640"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'4.90%"
$ws.Range("E25").Style = "Normal"
$ws.Range("F25").Value = "'6-2-2023"
$ws.Range("F25").Style = "Normal"
$ws.Range("G25").Value = "'0"
$ws.Range("G25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001229"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'3.03%"
$ws.Range("E26").Style = "Normal"
$ws.Range("F26").Value = "'6-2-2023"
$ws.Range("F26").Style = "Normal"
$ws.Range("G26").Value = "'0"
$ws.Range("G26").Style = "Normal"
$ws.Range("E27").Value = "'-0.42%"
$ws.Range("E27").Style = "Normal"
$ws.Range("F27").Value = "'6-2-2023"
$ws.Range("F27").Style = "Normal"
$ws.Range("G27").Value = "'0"
$ws.Range("G27").Style = "Normal"
$ws.Range("F28").Value = "'6-2-2023"
$ws.Range("F28").Style = "Normal"
$ws.Range("G28").Value = "'0"
$ws.Range("G28").Style = "Normal"
$ws.Range("F29").Value = "'6-2-2023"
$ws.Range("F29").Style = "Normal"
$ws.Range("G29").Value = "'0"
$ws.Range("G29").Style = "Normal"
$ws.Range("F30").Value = "'6-2-2023"
$ws.Range("F30").Style = "Normal"
$ws.Range("G30").Value = "'0"
$ws.Range("G30").Style = "Normal"
$ws.Range("F31").Value = "'6-2-2023"
$ws.Range("F31").Style = "Normal"
$ws.Range("G31").Value = "'0"
$ws.Range("G31").Style = "Normal"
$ws.Range("F32").Value = "'6-2-2023"
$ws.Range("F32").Style = "Normal"
$ws.Range("G32").Value = "'0"
$ws.Range("G32").Style = "Normal"
$ws.Range("F33").Value = "'6-2-2023"
$ws.Range("F33").Style = "Normal"
$ws.Range("G33").Value = "'0"
$ws.Range("G33").Style = "Normal"
$ws.Range("F34").Value = "'6-2-2023"
$ws.Range("F34").Style = "Normal"
$ws.Range("G34").Value = "'0"
$ws.Range("G34").Style = "Normal"
$ws.Range("F35").Value = "'6-2-2023"
$ws.Range("F35").Style = "Normal"
$ws.Range("G35").Value = "'0"
$ws.Range("G35").Style = "Normal"
$ws.Range("F36").Value = "'6-2-2023"
$ws.Range("F36").Style = "Normal"
$ws.Range("G36").Value = "'0"
$ws.Range("G36").Style = "Normal"
$ws.Range("F37").Value = "'6-2-2023"
$ws.Range("F37").Style = "Normal"
$ws.Range("G37").Value = "'0"
$ws.Range("G37").Style = "Normal"
$ws.Range("F38").Value = "'6-2-2023"
$ws.Range("F38").Style = "Normal"
$ws.Range("G38").Value = "'0"
$ws.Range("G38").Style = "Normal"
$ws.Range("D39").Value = "'0.02729"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-2.89%"
$ws.Range("E39").Style = "Normal"
$ws.Range("F39").Value = "'6-2-2023"
$ws.Range("F39").Style = "Normal"
$ws.Range("G39").Value = "'0"
$ws.Range("G39").Style = "Normal"
$ws.Range("D40").Value = "'0.05655"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'2.69%"
$ws.Range("E40").Style = "Normal"
$ws.Range("F40").Value = "'6-2-2023"
$ws.Range("F40").Style = "Normal"
$ws.Range("G40").Value = "'0"
$ws.Range("G40").Style = "Normal"
$ws.Range("E41").Value = "'25.96%"
$ws.Range("E41").Style = "Normal"
$ws.Range("F41").Value = "'6-2-2023"
$ws.Range("F41").Style = "Normal"
$ws.Range("G41").Value = "'0"
$ws.Range("G41").Style = "Normal"
$ws.Range("D42").Value = "'0.007687"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-3.32%"
$ws.Range("E42").Style = "Normal"
$ws.Range("F42").Value = "'6-2-2023"
$ws.Range("F42").Style = "Normal"
$ws.Range("G42").Value = "'0"
$ws.Range("G42").Style = "Normal"
$ws.Range("D43").Value = "'0.1403"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-1.83%"
$ws.Range("E43").Style = "Normal"
$ws.Range("F43").Value = "'6-2-2023"
$ws.Range("F43").Style = "Normal"
$ws.Range("G43").Value = "'0"
$ws.Range("G43").Style = "Normal"
$ws.Range("D44").Value = "'0.002013"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-6.74%"
$ws.Range("E44").Style = "Normal"
$ws.Range("F44").Value = "'6-2-2023"
$ws.Range("F44").Style = "Normal"
$ws.Range("G44").Value = "'0"
$ws.Range("G44").Style = "Normal"
$ws.Range("D45").Value = "'0.009431"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-6.72%"
$ws.Range("E45").Style = "Normal"
$ws.Range("F45").Value = "'6-2-2023"
$ws.Range("F45").Style = "Normal"
$ws.Range("G45").Value = "'0"
$ws.Range("G45").Style = "Normal"
$ws.Range("D46").Value = "'0.00007081"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-2.49%"
$ws.Range("E46").Style = "Normal"
$ws.Range("F46").Value = "'6-2-2023"
$ws.Range("F46").Style = "Normal"
$ws.Range("G46").Value = "'0"
$ws.Range("G46").Style = "Normal"
$ws.Range("E47").Value = "'-0.17%"
$ws.Range("E47").Style = "Normal"
$ws.Range("F47").Value = "'6-2-2023"
$ws.Range("F47").Style = "Normal"
$ws.Range("G47").Value = "'0"
$ws.Range("G47").Style = "Normal"
$ws.Range("D48").Value = "'0.003463"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-2.43%"
$ws.Range("E48").Style = "Normal"
$ws.Range("F48").Value = "'6-2-2023"
$ws.Range("F48").Style = "Normal"
$ws.Range("G48").Value = "'0"
$ws.Range("G48").Style = "Normal"
$ws.Range("E49").Value = "'-0.76%"
$ws.Range("E49").Style = "Normal"
$ws.Range("F49").Value = "'6-2-2023"
$ws.Range("F49").Style = "Normal"
$ws.Range("G49").Value = "'0"
$ws.Range("G49").Style = "Normal"
$ws.Range("E50").Value = "'-0.17%"
$ws.Range("E50").Style = "Normal"
$ws.Range("F50").Value = "'6-2-2023"
$ws.Range("F50").Style = "Normal"
$ws.Range("G50").Value = "'0"
$ws.Range("G50").Style = "Normal"
$ws.Range("E51").Value = "'-0.17%"
$ws.Range("E51").Style = "Normal"
$ws.Range("F51").Value = "'6-2-2023"
$ws.Range("F51").Style = "Normal"
$ws.Range("G51").Value = "'0"
$ws.Range("G51").Style = "Normal"
